# Update the "想去人数" (want-to-go count) figures that changed between
# scrapes, per the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 5743
$ws1.Range("F10").Value = 14
$ws1.Range("F13").Value = 1558
$ws1.Range("F14").Value = 1558
$ws1.Range("F18").Value = 584
$ws1.Range("F19").Value = 4151
$ws1.Range("F20").Value = 4151
$ws1.Range("F22").Value = 3316
$ws1.Range("F33").Value = 1118

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 244
$ws3.Range("F5").Value = 209

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 244
$ws4.Range("F8").Value  = 5743
$ws4.Range("F21").Value = 14
$ws4.Range("F23").Value = 1558
$ws4.Range("F24").Value = 1558
$ws4.Range("F29").Value = 584
$ws4.Range("F31").Value = 4151
$ws4.Range("F32").Value = 4151
$ws4.Range("F34").Value = 3316
$ws4.Range("F49").Value = 1118
